# OLX Monitor 2026-02-25 09:39 — append 7 new listing rows (227-233)
# to the "PODSUMOWANIE" detail log (columns A-H), extending the sheet
# from A1:H226 to A1:H233.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PODSUMOWANIE")

# --- 1. Stamp out the new rows with the same look & feel as the existing
#        log rows by copying the last populated row (226) down. This
#        carries the per-column styles (A=13 left, C=13 left, D/E=14
#        center, F=15 center/red "stale" highlight, B/G/H = no style).
$templateRow = $ws.Range("A226:H226")
$templateRow.Copy($ws.Range("A227:H227"))
$templateRow.Copy($ws.Range("A228:H228"))
$templateRow.Copy($ws.Range("A229:H229"))
$templateRow.Copy($ws.Range("A230:H230"))
$templateRow.Copy($ws.Range("A231:H231"))
$templateRow.Copy($ws.Range("A232:H232"))
$templateRow.Copy($ws.Range("A233:H233"))

# Row 230's "days listed" (F) is not flagged in the red/stale style
# (s=15) like the rest -- it uses the plain centered style (s=14, as
# seen e.g. on row 7) -- fix that cell's formatting before values go in.
$ws.Range("F7").Copy($ws.Range("F230"))

# --- 2. Fill in the actual scraped values for each new row.
$rows = @(
    @{ Row=227; A="2026-02-25 09:39:39"; B="poqui";           C="Kawalerka po remoncie z funkcjonalną antresolą - ul. Jana Sawy";                            D=2499;  E="28.10.2025"; F=119; G="https://www.olx.pl/d/oferta/kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger.html"; H="kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger" },
    @{ Row=228; A="2026-02-25 09:39:39"; B="poqui";           C="Przytulny pokój blisko Politechniki – ul. Przytulna";                                        D=549;   E="10.10.2025"; F=138; G="https://www.olx.pl/d/oferta/przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz.html"; H="przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz" },
    @{ Row=229; A="2026-02-25 09:39:39"; B="poqui";           C="Mieszkanie z KLIMATYZACJĄ 5 minut od UMCS, UP, KUL - Długosza";                               D=2049;  E="19.12.2025"; F=67;  G="https://www.olx.pl/d/oferta/mieszkanie-z-klimatyzacja-5-minut-od-umcs-up-kul-dlugosza-CID3-ID18KAEc.html"; H="mieszkanie-z-klimatyzacja-5-minut-od-umcs-up-kul-dlugosza-CID3-ID18KAEc" },
    @{ Row=230; A="2026-02-25 09:39:39"; B="pokojewlublinie"; C="WOLNY OD ZARAZ! Super lokalizacja, blisko centrum, ul. Paganiniego 12";                       D=12640; E="19.01.2026"; F=36;  G="https://www.olx.pl/d/oferta/wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc.html"; H="wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc" },
    @{ Row=231; A="2026-02-25 09:39:39"; B="pokojewlublinie"; C="WOLNY OD ZARAZ! Pokój jedynka, ul. Romanowskiego 58";                                        D=0;     E="11.08.2025"; F=197; G="https://www.olx.pl/d/oferta/wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm.html"; H="wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm" },
    @{ Row=232; A="2026-02-25 09:39:39"; B="dawnypatron";     C="Ładny pokój jednoosobowy. Wynajmę duży pokój w centrum. ul Niecała 4.";                      D=730;   E="20.09.2024"; F=522; G="https://www.olx.pl/d/oferta/ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM.html"; H="ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM" },
    @{ Row=233; A="2026-02-25 09:39:39"; B="dawnypatron";     C="Mam do wynajęcia pokój dla os. pracującej lub studenta. Narutowicza 14";                      D=14690; E="05.12.2025"; F=81;  G="https://www.olx.pl/d/oferta/mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv.html"; H="mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv" }
)

# Scratch cell (well outside the used A:H range) used as a bounce-pad so
# that "day.month.year" strings land as literal text, never as an
# auto-converted date serial. Excel's normal typed-value assignment
# (Range.Value = "...") runs the same smart-type heuristics a user
# typing into the grid would hit, which silently reinterprets
# unambiguous-looking strings like "10.10.2025" (day<=12, month<=12) as
# dates. Routing the literal through ="..."  (forces text result) and
# then PasteSpecial-ing *values only* into the real destination keeps
# the destination's existing number format/style (from the template
# copy above) while guaranteeing the stored value is plain text.
$scratch = $ws.Range("ZZ1")

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Cells.Item($n, 1).Value = $r.A
    $ws.Cells.Item($n, 2).Value = $r.B
    $ws.Cells.Item($n, 3).Value = $r.C
    $ws.Cells.Item($n, 4).Value = $r.D

    $scratch.Formula = '="' + $r.E + '"'
    $scratch.Copy()
    $ws.Cells.Item($n, 5).PasteSpecial(-4163)

    $ws.Cells.Item($n, 6).Value = $r.F
    $ws.Cells.Item($n, 7).Value = $r.G
    $ws.Cells.Item($n, 8).Value = $r.H
}

$scratch.ClearContents()
$excel.CutCopyMode = $false
